# Update "想去人数" (interest count) values on both the "展览" and
# "全部类型" worksheets to match the newly scraped data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F22").Value = 174
    $ws.Range("F26").Value = 4123
    $ws.Range("F33").Value = 495
    $ws.Range("F35").Value = 225
}
